$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing string for Thursday 16 maart
$ws.Range("D25").Value = "Ik was ziek"

# Add new rows for week of 20-23 maart
$ws.Range("A27").Value = "Maandag 20 maart: "
$ws.Range("D27").Value = "Absentie nog beter gemaakt "

$ws.Range("A29").Value = "Dinsdag 21 maart: "
$ws.Range("D29").Value = "Ervoor zorgen dat mensen je absentie kunnen afkeueren en goed keuren"

$ws.Range("A31").Value = "Woensdag 22 maart: "
$ws.Range("D31").Value = "Ervoor zorgen dat de absentie en het rooster beter worden gemaakt en dat het rooster op de goede dag komt van de week"

$ws.Range("A33").Value = "Donderdag 23 maart:"
$ws.Range("D33").Value = "Ervoor zorgen dat je de klas kan inzien en en daar weer het rooster van die klas kan zien en de leerlingen van de klas"

$ws.Range("D34").Select()
